$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New padron data replacing the old one (BAJA row removed, folios renumbered,
# "BAJA" entry gone entirely from the table).
$ws.Range("A1").Value = "Folio"
$ws.Range("B1").Value = "Nombre"

$ws.Range("A2").Value = 3564
$ws.Range("B2").Value = "CHOCOLATE"

$ws.Range("A3").Value = 3565
$ws.Range("B3").Value = "FINADO"

$ws.Range("A4").Value = 3566
$ws.Range("B4").Value = "INDEFINIDO"

$ws.Range("A5").Value = 3567
$ws.Range("B5").Value = "VAINILLA"

# Clear formatting from data rows (A2:B5) - no longer styled/bordered/filled
$ws.Range("A2:B5").ClearFormats()

# Row 6 becomes empty (last row previously held the removed "BAJA" entry)
$ws.Range("A6").ClearContents()
$ws.Range("B6").ClearContents()

# Update the selection to match the new state
$ws.Range("D10").Select()
